$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change D47:D68 completion status from "N" to "Y"
$ws.Range("D47:D68").Value = "Y"

# Update the view: scroll so row 57 is the top-left row, and select D72
$ws.Range("D72").Select()
$excel.ActiveWindow.ScrollRow = 57
